$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new worksheet "2022-Q3" right before the current "2022-Q2"
#    sheet (i.e. as the second tab, right after "总计"). All the existing
#    quarter sheets keep their data and simply shift one tab to the right.
# ---------------------------------------------------------------------------
$refSheet = $wb.Worksheets.Item(2)
$q3 = $wb.Worksheets.Add($refSheet)
$q3.Name = "2022-Q3"

# Header row (B1:H1) + the bold/bordered "index" column (A) use the same
# emphasised style as every other quarterly sheet in the workbook - build
# that look explicitly (bold, thin box border, centred/top aligned).
$header = $q3.Range("B1:H1")
$header.NumberFormat = "@"
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160
$header.Borders.LineStyle = 1

$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

$idxCol = $q3.Range("A2:A3")
$idxCol.Font.Bold = $true
$idxCol.HorizontalAlignment = -4108
$idxCol.VerticalAlignment = -4160
$idxCol.Borders.LineStyle = 1
$q3.Range("A2").Value = 0
$q3.Range("A3").Value = 1

# The fund-code / figure columns are stored as plain text in every quarterly
# sheet (so things like leading zeros and trailing zeros in "007581" /
# "0.5990" survive) - force text formatting before writing the values.
$body = $q3.Range("B2:G3")
$body.NumberFormat = "@"

$q3.Range("B2").Value = "213001"
$q3.Range("C2").Value = "宝盈鸿利收益灵活配置混合A"
$q3.Range("D2").Value = "16.06"
$q3.Range("E2").Value = "87.42"
$q3.Range("F2").Value = "3.73"
$q3.Range("G2").Value = "0.5990"
$q3.Range("H2").Value = 10

$q3.Range("B3").Value = "007581"
$q3.Range("C3").Value = "宝盈鸿利收益灵活配置混合C"
$q3.Range("D3").Value = "0.66"
$q3.Range("E3").Value = "87.42"
$q3.Range("F3").Value = "3.73"
$q3.Range("G3").Value = "0.0246"
$q3.Range("H3").Value = 10

# ---------------------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: shift the existing rows (2..8) down
#    by one row and insert the new 2022-Q3 summary figures into row 2.
#    Column A is just a running 0-based index (row-2), so it is rewritten
#    afterwards rather than shifted along with the rest of the row.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

for ($r = 8; $r -ge 2; $r--) {
    $src = $total.Range("B" + $r + ":D" + $r)
    $dst = $total.Range("B" + ($r + 1) + ":D" + ($r + 1))
    $src.Copy($dst)
}

$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.62

# Row 9 is brand new territory (previously the sheet only went to row 8), so
# give A9 the same style as the other column-A cells before writing values.
$total.Range("A8").Copy($total.Range("A9"))

for ($r = 2; $r -le 9; $r++) {
    $total.Range("A" + $r).Value = $r - 2
}
